$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: A15 now holds a label ("% local ... ") instead of being blank,
# and loses its bold formatting.
$ws.Range("A15").Value2 = '% local " " " " " '
$ws.Range("A15").Font.Bold = $false

# New summary rows 116-118 with additional percentage calculations.
$ws.Range("A116").Value2 = '% annon " " " " "'
$ws.Range("B116").Formula = "=(B106/SUM(B104:B113)) * 100"
$ws.Range("C116:P116").Formula = "=(C106/SUM(C104:C113)) * 100"

$ws.Range("A117").Value2 = '% local " " " " "'
$ws.Range("B117").Formula = "=(B105/SUM(B104:B113)) * 100"
$ws.Range("C117:P117").Formula = "=(C105/SUM(C104:C113)) * 100"

$ws.Range("A118").Value2 = '% nested + annon + local " " " " "'
$ws.Range("B118").Formula = "=SUM(B115:B117)"
$ws.Range("C118:P118").Formula = "=SUM(C115:C117)"

# Selection moves to B15 (single cell) with no frozen top-left anchor.
$ws.Range("B15").Select()
